$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 numeric values
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 92.0952380952381
$ws.Range("G2").Value = 105000
$ws.Range("H2").Value = 84000
$ws.Range("I2").Value = 21000
$ws.Range("N2").Value = 0.969

# Update row 2 time string values (also updates the shared strings table,
# which covers both E2/M2 pointing at the same text)
$ws.Range("E2").Value = "0:00:02"
$ws.Range("J2").Value = "0:08:41"
$ws.Range("K2").Value = "0:01:08"
$ws.Range("L2").Value = "0:00:07"
$ws.Range("M2").Value = "0:00:02"
